$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.047865867614746
$ws.Range("B1").Value = 4.405974864959717
$ws.Range("C1").Value = 3.825764894485474
$ws.Range("D1").Value = 1.725615382194519
$ws.Range("E1").Value = 0.9356040954589844
